$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3761.111
$ws.Range("J17").Value = 4132.9165
$ws.Range("L17").Value = 12398.7495
$ws.Range("N17").Value = -12734.7495

# ALC row 44
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 16000.0
$ws.Range("J44").Value = 16000.0
$ws.Range("L44").Value = 16000.0
$ws.Range("N44").Value = -16924.0

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 22224822.0
$ws.Range("I51").Value = 4499.5
$ws.Range("J51").Value = 37038372.0
$ws.Range("K51").Value = 4499.5
$ws.Range("L51").Value = 37038372.0
$ws.Range("M51").Value = -4015.5
$ws.Range("N51").Value = -37039340.0

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 487.22223
$ws.Range("I98").Value = 487.22223
$ws.Range("K98").Value = 487.22223
$ws.Range("M98").Value = 1010.77777

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2297.5
$ws.Range("I116").Value = 2211.6296
$ws.Range("J116").Value = 2761.2
$ws.Range("K116").Value = 2211.6296
$ws.Range("L116").Value = 2761.2
$ws.Range("M116").Value = 1230.3704
$ws.Range("N116").Value = -9645.2

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 487.22223
$ws.Range("I122").Value = 487.22223
$ws.Range("K122").Value = 1461.66669
$ws.Range("M122").Value = 988.33331

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 871.2963
$ws.Range("I129").Value = 312.85715
$ws.Range("J129").Value = 1066.75
$ws.Range("K129").Value = 938.5714499999999
$ws.Range("L129").Value = 3200.25
$ws.Range("M129").Value = 4061.42855
$ws.Range("N129").Value = -13200.25

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 730.3571
$ws.Range("I2").Value = 728.3333
$ws.Range("J2").Value = 742.5
$ws.Range("K2").Value = 728.3333
$ws.Range("L2").Value = 742.5
$ws.Range("M2").Value = -615.3333
$ws.Range("N2").Value = -968.5

# ARM row 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1758.5
$ws.Range("I21").Value = 1500.0
$ws.Range("K21").Value = 1500.0
$ws.Range("M21").Value = -1126.0

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21748.887
$ws.Range("I32").Value = 6017.6577
$ws.Range("J32").Value = 39863.637
$ws.Range("K32").Value = 6017.6577
$ws.Range("L32").Value = 39863.637
$ws.Range("M32").Value = -5730.6577
$ws.Range("N32").Value = -40437.637

# ARM row 56
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 0.0
$ws.Range("J56").Value = 0.0
$ws.Range("L56").Value = 0.0
$ws.Range("N56").ClearContents()

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1620.5454
$ws.Range("I110").Value = 1544.12
$ws.Range("J110").Value = 1859.375
$ws.Range("K110").Value = 1544.12
$ws.Range("L110").Value = 1859.375
$ws.Range("M110").Value = 500.8800000000001
$ws.Range("N110").Value = -5949.375

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 730.3571
$ws.Range("I116").Value = 728.3333
$ws.Range("J116").Value = 742.5
$ws.Range("K116").Value = 728.3333
$ws.Range("L116").Value = 742.5
$ws.Range("M116").Value = 1565.6667
$ws.Range("N116").Value = -5330.5

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 730.3571
$ws.Range("I3").Value = 728.3333
$ws.Range("J3").Value = 742.5
$ws.Range("K3").Value = 728.3333
$ws.Range("L3").Value = 742.5
$ws.Range("M3").Value = -614.3333
$ws.Range("N3").Value = -970.5

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3736.7036
$ws.Range("I99").Value = 3834.55
$ws.Range("J99").Value = 3457.1428
$ws.Range("K99").Value = 3834.55
$ws.Range("L99").Value = 3457.1428
$ws.Range("M99").Value = -2336.55
$ws.Range("N99").Value = -6453.1428

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3736.7036
$ws.Range("I126").Value = 3834.55
$ws.Range("J126").Value = 3457.1428
$ws.Range("K126").Value = 11503.65
$ws.Range("L126").Value = 10371.4284
$ws.Range("M126").Value = -9033.650000000001
$ws.Range("N126").Value = -15311.4284

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1940.2333
$ws.Range("I132").Value = 1330.2941
$ws.Range("K132").Value = 3990.8823
$ws.Range("M132").Value = -1460.8823

# CUL row 62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4757.0
$ws.Range("J62").Value = 4757.0
$ws.Range("L62").Value = 14271.0
$ws.Range("N62").Value = -15643.0

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 9562.8
$ws.Range("I63").Value = 0.0
$ws.Range("J63").Value = 9562.8
$ws.Range("K63").Value = 0.0
$ws.Range("L63").Value = 28688.4
$ws.Range("N63").Value = -30186.4
$ws.Range("M63").ClearContents()

# CUL row 65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 4757.0
$ws.Range("J65").Value = 4757.0
$ws.Range("L65").Value = 42813.0
$ws.Range("N65").Value = -49677.0

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 9562.8
$ws.Range("I66").Value = 0.0
$ws.Range("J66").Value = 9562.8
$ws.Range("K66").Value = 0.0
$ws.Range("L66").Value = 86065.2
$ws.Range("N66").Value = -93553.2
$ws.Range("M66").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3424.25
$ws.Range("I80").Value = 5043.5713
$ws.Range("J80").Value = 2552.3076
$ws.Range("K80").Value = 5043.5713
$ws.Range("L80").Value = 2552.3076
$ws.Range("M80").Value = -4045.5713
$ws.Range("N80").Value = -4548.3076

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3424.25
$ws.Range("I83").Value = 5043.5713
$ws.Range("J83").Value = 2552.3076
$ws.Range("K83").Value = 25217.8565
$ws.Range("L83").Value = 12761.538
$ws.Range("M83").Value = -20225.8565
$ws.Range("N83").Value = -22745.538

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2417.8333
$ws.Range("I122").Value = 3076.75
$ws.Range("J122").Value = 1100.0
$ws.Range("K122").Value = 9230.25
$ws.Range("L122").Value = 3300.0
$ws.Range("M122").Value = -6780.25
$ws.Range("N122").Value = -8200.0

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2104.6875
$ws.Range("I126").Value = 1853.4286
$ws.Range("J126").Value = 2584.3635
$ws.Range("K126").Value = 5560.2858
$ws.Range("L126").Value = 7753.0905
$ws.Range("M126").Value = -3090.2858
$ws.Range("N126").Value = -12693.0905

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2114.0
$ws.Range("I7").Value = 1168.0
$ws.Range("K7").Value = 1168.0
$ws.Range("M7").Value = -1056.0

# LTW row 39
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 14982.5
$ws.Range("J39").Value = 14982.5
$ws.Range("L39").Value = 14982.5
$ws.Range("N39").Value = -15902.5

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2318.3333
$ws.Range("I82").Value = 1555.0
$ws.Range("J82").Value = 2700.0
$ws.Range("K82").Value = 1555.0
$ws.Range("L82").Value = 2700.0
$ws.Range("M82").Value = -1194.0
$ws.Range("N82").Value = -3422.0

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2318.3333
$ws.Range("I85").Value = 1555.0
$ws.Range("J85").Value = 2700.0
$ws.Range("K85").Value = 1555.0
$ws.Range("L85").Value = 2700.0
$ws.Range("M85").Value = -307.0
$ws.Range("N85").Value = -5196.0

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2114.0
$ws.Range("I126").Value = 1168.0
$ws.Range("K126").Value = 3504.0
$ws.Range("M126").Value = -1034.0

# WVR row 15
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8005.3335
$ws.Range("J15").Value = 8005.3335
$ws.Range("L15").Value = 8005.3335
$ws.Range("N15").Value = -8581.3335

# WVR row 23
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 703.3333
$ws.Range("I23").Value = 703.3333
$ws.Range("K23").Value = 703.3333
$ws.Range("M23").Value = -474.3333
